# fix: Update group type and classe to not null
#
# The "group" worksheet is missing "type" and "classe" columns (they were
# null/absent). This adds the two columns with values for every existing
# data row, then leaves the "group" sheet active/selected as the last
# sheet the (human) editor was working on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group")

# Header row (row 1)
$ws.Cells.Item(1, 5).Value2 = "type"
$ws.Cells.Item(1, 6).Value2 = "classe"

# Data rows (rows 2-5): every group is a "CM" (cours magistral) in class "GB"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value2 = "CM"
    $ws.Cells.Item($r, 6).Value2 = "GB"
}

# The editor ended up with the "group" sheet active, cell H3 selected.
$ws.Activate()
$ws.Range("H3").Select()
